# Updating changes related to env setup
# Restructures the "prodfix" sheet's Extraction-file/Reported-Variables
# block: renames a couple of headers and introduces a new
# "Expected_File_names" column (with the bare file name, derived from the
# full path already present in the "Files_to_upload" column) between the
# existing "Files_to_upload" and "ReportedVariables" columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("prodfix")

# --- Column B: rename header + update the repeated value in data rows ---
$ws.Range("B1").Value = "Population_name"
$ws.Range("B2").Value = "PRODFix_QOL_ECON - UtilityOutcome - 9/19/2022"
$ws.Range("B5").Value = "PRODFix_QOL_ECON - UtilityOutcome - 9/19/2022"
$ws.Range("B8").Value = "PRODFix_QOL_ECON - UtilityOutcome - 9/19/2022"
$ws.Range("B11").Value = "PRODFix_QOL_ECON - UtilityOutcome - 9/19/2022"

# --- Column H: rename header (column itself does not move) ---
$ws.Range("H1").Value = "Files_to_upload"

# --- Insert a new column I ("Expected_File_names"), shifting old I:K -> J:L ---
$ws.Columns("I").Insert()
$ws.Columns("I").ColumnWidth = 36.6

$ws.Range("I1").Value = "Expected_File_names"
$ws.Range("I2").Value = "UtilityOutcome_Feature_Extraction_file_QoL_UtilityData_ECON_NoUtility.xlsx"
$ws.Range("I5").Value = "UtilityOutcome_Feature_Extraction_file_ECON_UtilityData_QoL_NoUtility.xlsx"
$ws.Range("I8").Value = "UtilityOutcome_Feature_Extraction_file_Both_QoL_ECON_Utility.xlsx"
$ws.Range("I11").Value = "UtilityOutcome_Feature_Extraction_file_NegativeScenario_QoL__ECON_NoUtility.xlsx"

# --- Update the view state: scrolled right toward column G, selection on I11 ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 7
$ws.Range("I11").Select()
